$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false,
                             $true, 1, $false, $replace, 2)
}

Replace-Text "2026-01-05 Monday" "2026-01-06 Tuesday"

Replace-Text "34×35=" "62×91="
Replace-Text "28×15=" "98×53="
Replace-Text "84×97=" "83×62="
Replace-Text "76×25=" "35×73="
Replace-Text "11×25=" "42×44="
Replace-Text "40×56=" "41×12="
Replace-Text "87×50=" "26×19="
Replace-Text "12×67=" "91×63="
Replace-Text "11×21=" "40×96="
Replace-Text "81×62=" "21×96="
Replace-Text "79×81=" "50×25="
Replace-Text "51×96=" "46×36="
Replace-Text "19×68=" "83×19="
Replace-Text "13×13=" "23×62="
Replace-Text "15×61=" "61×87="
Replace-Text "77×27=" "89×60="
Replace-Text "36×82=" "51×92="
Replace-Text "76×42=" "18×99="
Replace-Text "82×88=" "77×99="
Replace-Text "25×75=" "85×80="
Replace-Text "24×58=" "40×76="
Replace-Text "24×56=" "44×65="
Replace-Text "27×12=" "16×27="
Replace-Text "25×74=" "65×88="
Replace-Text "33×25=" "82×18="
